# ------------------------------------------------------------------
# Applies the commit's changes to PlayerPerformance_3829.xlsx:
#   1. "ODI Batting": drop the placeholder (empty) INNING_NUMBER cells
#      in column B that never held a value, and normalize the stray
#      non-breaking space in E5 to a regular space.
#   2. Add a new "ODI Batting Extra" worksheet (after "ODI Bowling")
#      with per-match extra batting stats.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Tidy up the existing "ODI Batting" sheet -------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")

# These B-column cells are blank placeholders (no INNING_NUMBER
# recorded) - clear them so they no longer appear as empty cell
# nodes in the saved sheet.
$emptyInningCells = @("B4","B5","B7","B9","B10","B12","B15","B16","B17","B18","B19","B20")
foreach ($cellRef in $emptyInningCells) {
    $battingSheet.Range($cellRef).ClearContents()
}

# E5 (MATCH_INNING) held a stray non-breaking space - normalize it to
# a plain space character.
$battingSheet.Range("E5").Value = " "

# --- 2. Add the new "ODI Batting Extra" worksheet -------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extraSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$extraSheet.Name = "ODI Batting Extra"

# Header row (row 1)
$headers = @("MATCH_CODE","BATTING_POSITION","NUM_4","NUM_6","PERCENT_RUNS_OF_TOTAL","MAN_OF_MATCH")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $extraSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Match the header formatting used on the other sheets (bold, thin
# border, centered) by copying the style from an existing header cell.
$battingSheet.Range("A1").Copy()
$extraSheet.Range("A1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Helper: write a value as literal text so Excel doesn't silently
# convert numeric-looking strings (e.g. "3188", "0", "0.96%") into
# numbers / percentages.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# MATCH_CODE (column A) for every row, taken from the ODI Batting sheet
$matchCodes = @("3188","3190","3220","3223","3225","3228","3230","3232","3236","3242","3252","3272","3607","3609","4024","4027","4402","4406","4410")
for ($i = 0; $i -lt $matchCodes.Length; $i++) {
    Set-TextValue $extraSheet.Cells.Item($i + 2, 1) $matchCodes[$i]
}

# MAN_OF_MATCH (column F) is known ("NO") for the first 8 matches only
for ($row = 2; $row -le 9; $row++) {
    Set-TextValue $extraSheet.Cells.Item($row, 6) "NO"
}

# BATTING_POSITION (column B) - numeric, recorded for rows 4-6
$extraSheet.Range("B4").Value = 11
$extraSheet.Range("B5").Value = 11
$extraSheet.Range("B6").Value = 11

# NUM_4 / NUM_6 (columns C & D) - recorded for rows 5-6
Set-TextValue $extraSheet.Range("C5") "0"
Set-TextValue $extraSheet.Range("D5") "0"
Set-TextValue $extraSheet.Range("C6") "0"
Set-TextValue $extraSheet.Range("D6") "0"

# PERCENT_RUNS_OF_TOTAL (column E) - recorded for row 6 only
Set-TextValue $extraSheet.Range("E6") "0.96%"
